$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("M2").Value = 0.3737363333333333
$ws.Range("N2").Value = 1.121209
$ws.Range("O2").Value = 0.0566058187608514
$ws.Range("P2").Value = 0.05660581876085141
$ws.Range("Q2").Value = 23.88593094425811
$ws.Range("R2").Value = 214.973378498323
$ws.Range("S2").Value = 0.02302683334243916
$ws.Range("T2").Value = 0.02302683334243917

$ws.Range("G3").Value = 63.91118233333333
$ws.Range("H3").Value = 191.733547
$ws.Range("I3").Value = 0.4067926910433548
$ws.Range("J3").Value = 0.4067926910433549
$ws.Range("O3").Value = 0.0001561544702435616
$ws.Range("P3").Value = 0.0001561544702435616
$ws.Range("Q3").Value = 0.06589242898566666
$ws.Range("R3").Value = 0.5930318608709999
$ws.Range("S3").Value = 0.00006352249716882788
$ws.Range("T3").Value = 0.0000635224971688279

$ws.Range("G4").Value = 63.91118233333333
$ws.Range("H4").Value = 191.733547
$ws.Range("I4").Value = 0.4067926910433548
$ws.Range("J4").Value = 0.4067926910433549
$ws.Range("M4").Value = 6.227669333333334
$ws.Range("N4").Value = 18.683008
$ws.Range("O4").Value = 0.943238026768905
$ws.Range("P4").Value = 0.943238026768905
$ws.Range("Q4").Value = 398.0177102743751
$ws.Range("R4").Value = 3582.159392469376
$ws.Range("S4").Value = 0.3837023352037468
$ws.Range("T4").Value = 0.3837023352037469

$ws.Range("I5").Value = 0.3656254573230189
$ws.Range("J5").Value = 0.365625457323019
$ws.Range("M5").Value = 0.3737363333333333
$ws.Range("N5").Value = 1.121209
$ws.Range("O5").Value = 0.0566058187608514
$ws.Range("P5").Value = 0.05660581876085141
$ws.Range("Q5").Value = 21.4686856902
$ws.Range("R5").Value = 193.2181712118
$ws.Range("S5").Value = 0.02069652837158022
$ws.Range("T5").Value = 0.02069652837158023

$ws.Range("I6").Value = 0.3656254573230189
$ws.Range("J6").Value = 0.365625457323019
$ws.Range("O6").Value = 0.0001561544702435616
$ws.Range("P6").Value = 0.0001561544702435616
$ws.Range("R6").Value = 0.5330173086
$ws.Range("S6").Value = 0.00005709404959583596
$ws.Range("T6").Value = 0.00005709404959583596

$ws.Range("I7").Value = 0.3656254573230189
$ws.Range("J7").Value = 0.365625457323019
$ws.Range("M7").Value = 6.227669333333334
$ws.Range("N7").Value = 18.683008
$ws.Range("O7").Value = 0.943238026768905
$ws.Range("P7").Value = 0.943238026768905
$ws.Range("Q7").Value = 357.7385005824
$ws.Range("R7").Value = 3219.6465052416
$ws.Range("S7").Value = 0.3448718349018429
$ws.Range("T7").Value = 0.3448718349018429

$ws.Range("G8").Value = 35.755375
$ws.Range("H8").Value = 107.266125
$ws.Range("I8").Value = 0.2275818516336261
$ws.Range("J8").Value = 0.2275818516336262
$ws.Range("M8").Value = 0.3737363333333333
$ws.Range("N8").Value = 1.121209
$ws.Range("O8").Value = 0.0566058187608514
$ws.Range("P8").Value = 0.05660581876085141
$ws.Range("Q8").Value = 13.36308274945833
$ws.Range("R8").Value = 120.267744745125
$ws.Range("S8").Value = 0.01288245704683201
$ws.Range("T8").Value = 0.01288245704683202

$ws.Range("G9").Value = 35.755375
$ws.Range("H9").Value = 107.266125
$ws.Range("I9").Value = 0.2275818516336261
$ws.Range("J9").Value = 0.2275818516336262
$ws.Range("O9").Value = 0.0001561544702435616
$ws.Range("P9").Value = 0.0001561544702435616
$ws.Range("Q9").Value = 0.036863791625
$ws.Range("R9").Value = 0.331774124625
$ws.Range("S9").Value = 0.00003553792347889772
$ws.Range("T9").Value = 0.00003553792347889772

$ws.Range("G10").Value = 35.755375
$ws.Range("H10").Value = 107.266125
$ws.Range("I10").Value = 0.2275818516336261
$ws.Range("J10").Value = 0.2275818516336262
$ws.Range("M10").Value = 6.227669333333334
$ws.Range("N10").Value = 18.683008
$ws.Range("O10").Value = 0.943238026768905
$ws.Range("P10").Value = 0.943238026768905
$ws.Range("Q10").Value = 222.6726523893334
$ws.Range("R10").Value = 2004.053871504
$ws.Range("S10").Value = 0.2146638566633152
$ws.Range("T10").Value = 0.2146638566633152
